# Work diary for the last 45 minutes: append a new entry (row) describing
# error/exception handling work to the "Tableau1" table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Tableau1")

# Grow the table by one row - this also extends the table ref, the
# autofilter range and the sheet dimension, same as typing into the first
# blank row below an existing table in Excel.
$newRow = $table.ListRows.Add()

# Bring over the date-column formatting (numFmtId 14 + wrap text) from the
# row above before putting any other new formatting in play, so the new
# cell reuses the existing style instead of Excel minting a fresh one.
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)

$ws.Range("A49").Value = 44694
$ws.Range("B49").Value = "Réalisation"
$ws.Range("C49").Value = 0.75
$ws.Range("D49").Value = "Gestions des erreurs/exceptions"
$ws.Range("E49").Value = "Affichage du message d'erreur a l'utilisateur"

$ws.Range("F49").Value = "76e90026b3f5ee849f3a2f6bb866614cdf21086e"
$ws.Range("F49").NumberFormat = "0.00E+00"
$ws.Range("F49").WrapText = $true

$ws.Range("A50").Select()
